$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to reduced-precision (custom accuracy) figures
$ws.Range("B5").Value = 12.97
$ws.Range("C5").Value = 9.4
$ws.Range("E5").Value = 28.13
$ws.Range("F5").Value = 22.92
$ws.Range("G5").Value = 10.21
$ws.Range("I5").Value = 15.71
$ws.Range("J5").Value = 6.91
$ws.Range("K5").Value = 10.17
$ws.Range("L5").Value = 11.3
$ws.Range("M5").Value = 11.85
$ws.Range("N5").Value = 3.26
$ws.Range("O5").Value = 10.15
$ws.Range("P5").Value = 14.4
$ws.Range("Q5").Value = 8.67
$ws.Range("R5").Value = 0.79
$ws.Range("S5").Value = 0.61
$ws.Range("T5").Value = 147.2
$ws.Range("U5").Value = 28.49
$ws.Range("W5").Value = 19.03
$ws.Range("X5").Value = 9.949999999999999
$ws.Range("Y5").Value = 1.63
$ws.Range("Z5").Value = 20.24
$ws.Range("AA5").Value = 8.279999999999999
$ws.Range("AB5").Value = 7.41
$ws.Range("AD5").Value = 11.79
$ws.Range("AE5").Value = 0.5600000000000001
$ws.Range("AF5").Value = 38.94
$ws.Range("AG5").Value = 5.23
$ws.Range("AH5").Value = 11.71

# Remove row 6 entirely (data trimmed from 2 rows to 1 row of readings)
$ws.Rows(6).Delete()
